$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "donor" sample row (k@gmail.com / Ken / Glin) is removed; the real
# upload rows (John Doe, Kayla Coms) shift up and get refreshed
# Amount/Date values so the sheet now reflects real upload data.
$ws.Rows(2).Delete()

# Row 2 is now John Doe -- refresh Amount and Date.
$ws.Range("D2").Value = 10
$ws.Range("E2").Value = "2021-04-24 15:49:41"

# Row 3 is now Kayla Coms -- refresh Amount and Date.
$ws.Range("D3").Value = 5
$ws.Range("E3").Value = "2021-04-12 11:49:41"

# Rebuild the hyperlinks collection so it only covers the two remaining
# rows (dropping the stale entry that pointed at the now-removed row 4).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:k@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:jd@gmail.com") | Out-Null

# Re-adding the hyperlinks resets number formatting on those cells; put it
# back the way it was so A2/A3 keep their original text format.
$ws.Range("A2:A3").NumberFormat = "@"

# Match the saved selection state.
$ws.Range("C4").Select()
